$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 664.0833
$ws.Range("I98").Value = 664.0833
$ws.Range("K98").Value = 664.0833
$ws.Range("M98").Value = 833.9167
$ws.Range("H106").Value = 3100
$ws.Range("I106").Value = 3100
$ws.Range("K106").Value = 3100
$ws.Range("M106").Value = -2469
$ws.Range("H116").Value = 5248.75
$ws.Range("I116").Value = 4998.5
$ws.Range("K116").Value = 4998.5
$ws.Range("M116").Value = -1556.5
$ws.Range("H122").Value = 664.0833
$ws.Range("I122").Value = 664.0833
$ws.Range("K122").Value = 1992.2499
$ws.Range("M122").Value = 457.7501
$ws.Range("H125").Value = 4361.2856
$ws.Range("I125").Value = 532
$ws.Range("J125").Value = 4999.5
$ws.Range("K125").Value = 4788
$ws.Range("L125").Value = 44995.5
$ws.Range("M125").Value = -2328
$ws.Range("N125").Value = -49915.5
$ws.Range("H132").Value = 1186.0869
$ws.Range("I132").Value = 1103.6364
$ws.Range("K132").Value = 3310.9092
$ws.Range("M132").Value = -780.9092000000001
$ws.Range("H135").Value = 1020.0909
$ws.Range("I135").Value = 1023
$ws.Range("J135").Value = 1007
$ws.Range("K135").Value = 9207
$ws.Range("L135").Value = 9063
$ws.Range("M135").Value = -6672
$ws.Range("N135").Value = -14133
$ws.Range("H138").Value = 2720.1667
$ws.Range("I138").Value = 1112.1786
$ws.Range("K138").Value = 3336.5358
$ws.Range("M138").Value = 1803.4642

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 95.5
$ws.Range("I5").Value = 95.5
$ws.Range("K5").Value = 95.5
$ws.Range("M5").Value = 16.5
$ws.Range("H32").Value = 4917.2334
$ws.Range("I32").Value = 4197.4287
$ws.Range("K32").Value = 4197.4287
$ws.Range("M32").Value = -3910.4287
$ws.Range("H97").Value = 825.25
$ws.Range("I97").Value = 420.07693
$ws.Range("J97").Value = 2581
$ws.Range("K97").Value = 420.07693
$ws.Range("L97").Value = 2581
$ws.Range("M97").Value = 75.92307
$ws.Range("N97").Value = -3573

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 95.5
$ws.Range("I4").Value = 95.5
$ws.Range("K4").Value = 95.5
$ws.Range("M4").Value = 19.5
$ws.Range("H20").Value = 2854.8125
$ws.Range("I20").Value = 2790.2727
$ws.Range("K20").Value = 2790.2727
$ws.Range("M20").Value = -2543.2727
$ws.Range("H94").Value = 366.42856
$ws.Range("I94").Value = 398
$ws.Range("J94").Value = 287.5
$ws.Range("K94").Value = 398
$ws.Range("L94").Value = 287.5
$ws.Range("M94").Value = 53
$ws.Range("N94").Value = -1189.5
$ws.Range("H99").Value = 699.6667
$ws.Range("I99").Value = 500
$ws.Range("J99").Value = 799.5
$ws.Range("K99").Value = 500
$ws.Range("L99").Value = 799.5
$ws.Range("M99").Value = 998
$ws.Range("N99").Value = -3795.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 659.8333
$ws.Range("I22").Value = 972
$ws.Range("J22").Value = 597.4
$ws.Range("K22").Value = 972
$ws.Range("L22").Value = 597.4
$ws.Range("M22").Value = -622
$ws.Range("N22").Value = -1297.4
$ws.Range("H31").Value = 2617.875
$ws.Range("I31").Value = 2707.7144
$ws.Range("K31").Value = 2707.7144
$ws.Range("M31").Value = -2412.7144
$ws.Range("H34").Value = 2617.875
$ws.Range("I34").Value = 2707.7144
$ws.Range("K34").Value = 2707.7144
$ws.Range("M34").Value = -2505.7144
$ws.Range("H86").Value = 17369.95
$ws.Range("I86").Value = 4514.6665
$ws.Range("J86").Value = 36652.875
$ws.Range("K86").Value = 4514.6665
$ws.Range("L86").Value = 36652.875
$ws.Range("M86").Value = -3391.6665
$ws.Range("N86").Value = -38898.875
$ws.Range("H89").Value = 17369.95
$ws.Range("I89").Value = 4514.6665
$ws.Range("J89").Value = 36652.875
$ws.Range("K89").Value = 22573.3325
$ws.Range("L89").Value = 183264.375
$ws.Range("M89").Value = -16957.3325
$ws.Range("N89").Value = -194496.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 695.4286
$ws.Range("I5").Value = 644.6667
$ws.Range("K5").Value = 1934.0001
$ws.Range("M5").Value = -1822.0001
$ws.Range("H8").Value = 460
$ws.Range("I8").Value = 460
$ws.Range("K8").Value = 1380
$ws.Range("M8").Value = -1241
$ws.Range("H12").Value = 154.6875
$ws.Range("I12").Value = 181.85715
$ws.Range("J12").Value = 133.55556
$ws.Range("K12").Value = 545.5714499999999
$ws.Range("L12").Value = 400.66668
$ws.Range("M12").Value = -372.5714499999999
$ws.Range("N12").Value = -746.66668
$ws.Range("H23").Value = 323.83334
$ws.Range("J23").Value = 328.6
$ws.Range("L23").Value = 985.8000000000001
$ws.Range("N23").Value = -1455.8
$ws.Range("H33").Value = 740.875
$ws.Range("I33").Value = 132.14285
$ws.Range("J33").Value = 5002
$ws.Range("K33").Value = 792.8571000000001
$ws.Range("L33").Value = 30012
$ws.Range("M33").Value = -509.8571000000001
$ws.Range("N33").Value = -30578
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9470
$ws.Range("H68").Value = 1200
$ws.Range("J68").Value = 1200
$ws.Range("L68").Value = 3600
$ws.Range("N68").Value = -5222
$ws.Range("H71").Value = 1200
$ws.Range("J71").Value = 1200
$ws.Range("L71").Value = 10800
$ws.Range("N71").Value = -18912
$ws.Range("H86").Value = 292.25
$ws.Range("J86").Value = 285
$ws.Range("L86").Value = 855
$ws.Range("N86").Value = -3227
$ws.Range("H89").Value = 292.25
$ws.Range("J89").Value = 285
$ws.Range("L89").Value = 2565
$ws.Range("N89").Value = -14421
$ws.Range("H122").Value = 802.3333
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 9000
$ws.Range("N122").Value = -13900
$ws.Range("H125").Value = 1999.5
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H135").Value = 695.4286
$ws.Range("I135").Value = 644.6667
$ws.Range("K135").Value = 5802.0003
$ws.Range("M135").Value = -3267.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5266.6665
$ws.Range("I102").Value = 5266.6665
$ws.Range("K102").Value = 5266.6665
$ws.Range("M102").Value = -3644.6665
$ws.Range("H113").Value = 1749.5
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 1499
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 1499
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -5839
$ws.Range("H122").Value = 776.5
$ws.Range("I122").Value = 776.5
$ws.Range("K122").Value = 2329.5
$ws.Range("M122").Value = 120.5
$ws.Range("H126").Value = 2036.3334
$ws.Range("I126").Value = 2036.3334
$ws.Range("K126").Value = 6109.0002
$ws.Range("M126").Value = -3639.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4747.5
$ws.Range("I40").Value = 4747.5
$ws.Range("K40").Value = 4747.5
$ws.Range("M40").Value = -4611.5
$ws.Range("H122").Value = 2826.7273
$ws.Range("I122").Value = 2344.889
$ws.Range("K122").Value = 7034.667
$ws.Range("M122").Value = -4584.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 999.5
$ws.Range("I107").Value = 999.5
$ws.Range("K107").Value = 2998.5
$ws.Range("M107").Value = -1078.5
$ws.Range("H126").Value = 2627.6667
$ws.Range("I126").Value = 2627.6667
$ws.Range("K126").Value = 7883.000100000001
$ws.Range("M126").Value = -5413.000100000001
$ws.Range("H132").Value = 1293.7727
$ws.Range("I132").Value = 1261.9445
$ws.Range("J132").Value = 1437
$ws.Range("K132").Value = 3785.8335
$ws.Range("L132").Value = 4311
$ws.Range("M132").Value = -1255.8335
$ws.Range("N132").Value = -9371
